{"js": "// Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table.\n// The table has 20 rows \u00d7 5 columns; only rows 0, 4, 8, 12, 16 (0-based)\n// hold text (5 values each = 25 cells total). Each cell's text is replaced,\n// in document (row-major) order, with its new value.\n\nconst newValues = [\n  \"63\u00f75=12, 3\", \"63\u00f79=7, 0\", \"66\u00f77=9, 3\", \"89\u00f77=12, 5\", \"60\u00f73=20, 0\",\n  \"32\u00f72=16, 0\", \"66\u00f77=9, 3\", \"89\u00f75=17, 4\", \"19\u00f76=3, 1\", \"82\u00f79=9, 1\",\n  \"80\u00f78=10, 0\", \"23\u00f75=4, 3\", \"61\u00f73=20, 1\", \"94\u00f73=31, 1\", \"65\u00f73=21, 2\",\n  \"89\u00f78=11, 1\", \"49\u00f72=24, 1\", \"51\u00f77=7, 2\", \"82\u00f79=9, 1\", \"75\u00f77=10, 5\",\n  \"64\u00f78=8, 0\", \"56\u00f79=6, 2\", \"18\u00f72=9, 0\", \"53\u00f79=5, 8\", \"13\u00f77=1, 6\",\n];\n\nconst rowIndexes = [0, 4, 8, 12, 16];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nlet k = 0;\nfor (const rowIdx of rowIndexes) {\n  for (let col = 0; col < 5; col++) {\n    table.getCell(rowIdx, col).value = newValues[k];\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"two-digit / one-digit\" answer cells in the single table.\n# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-based)\n# hold text (5 values each = 25 cells total). Each cell's text is replaced\n# in document (row-major) order with its new value.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$newValues = @(\n  \"63\u00f75=12, 3\", \"63\u00f79=7, 0\", \"66\u00f77=9, 3\", \"89\u00f77=12, 5\", \"60\u00f73=20, 0\",\n  \"32\u00f72=16, 0\", \"66\u00f77=9, 3\", \"89\u00f75=17, 4\", \"19\u00f76=3, 1\", \"82\u00f79=9, 1\",\n  \"80\u00f78=10, 0\", \"23\u00f75=4, 3\", \"61\u00f73=20, 1\", \"94\u00f73=31, 1\", \"65\u00f73=21, 2\",\n  \"89\u00f78=11, 1\", \"49\u00f72=24, 1\", \"51\u00f77=7, 2\", \"82\u00f79=9, 1\", \"75\u00f77=10, 5\",\n  \"64\u00f78=8, 0\", \"56\u00f79=6, 2\", \"18\u00f72=9, 0\", \"53\u00f79=5, 8\", \"13\u00f77=1, 6\"\n)\n\n$rowIndexes = @(1, 5, 9, 13, 17)\n\n$k = 0\nforeach ($r in $rowIndexes) {\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$k]\n    $k = $k + 1\n  }\n}\n"}
